$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in column D values that were previously blank/missing
$ws.Range("D3").Value = 15809
$ws.Range("D4").Value = 10186
$ws.Range("D5").Value = 9163
$ws.Range("D6").Value = 161
$ws.Range("D8").Value = 3621
$ws.Range("D9").Value = 3144
$ws.Range("D10").Value = 994
$ws.Range("D11").Value = 500
$ws.Range("D13").Value = 2903
$ws.Range("D14").Value = 2932
$ws.Range("D15").Value = 1547
$ws.Range("D16").Value = 224
$ws.Range("D18").Value = 183561
$ws.Range("D19").Value = 104
$ws.Range("D20").Value = 194591
$ws.Range("D21").Value = 184538
$ws.Range("D22").Value = 1752
$ws.Range("D23").Value = 1729
$ws.Range("D24").Value = 265
$ws.Range("D26").Value = 4980
$ws.Range("D27").Value = 100000000
$ws.Range("D28").Value = 536
$ws.Range("D30").Value = 3830
$ws.Range("D31").Value = 271
$ws.Range("D33").Value = 5061
$ws.Range("D34").Value = 122

# Update the sheet view: scroll to show row 25 at top and select D35
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
